$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Login" test data: several rows' Actual Status column (D)
# now reflect a different validation message returned by the API.
$ws.Range("D10").Value = "Invalid parameters provided. Please provide valid parameters."
$ws.Range("D11").Value = "Invalid parameters provided. Please provide valid parameters."
$ws.Range("D12").Value = "Invalid parameters provided. Please provide valid parameters."
$ws.Range("D16").Value = "Invalid Email Address "
